# B6-PowerPoint.pptx — Fri, May 15, 2020  8:05:30 PM
#
# The three data tables (on slides 14, 15 and 16) that were using the
# deck's custom "Table_0" table style ({5EFFE562-3217-48EB-B695-77BA6CED0536})
# get switched over to a different table style
# ({8AA42090-D2FB-4EF7-996E-B46A52E4E2CE}).
#
# PowerPoint doesn't allow Table.Style to be assigned directly (it raises
# "Table styles cannot be assigned through a property" if you try) — the
# supported COM call is Table.ApplyStyle("{GUID}").

$p = $ppt.ActivePresentation

$targetStyleId = "{8AA42090-D2FB-4EF7-996E-B46A52E4E2CE}"
$slideIndexes = @(14, 15, 16)

foreach ($slideIdx in $slideIndexes) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
